# Continue the "Mitschrift" CSS notes section with the additional bullet
# points captured at the end of the lecture (margin/padding box model,
# shared-class advice, borders, border-radius) and relocate the trailing
# "_GoBack" bookmark so it still marks the very end of the document.

$d = $word.ActiveDocument

# Locate the insertion point: right after "Padding: " and before the
# trailing bookmark that currently sits at the end of the document.
$find = $d.Content
[void]$find.Find.Execute("Padding: ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$find.Collapse(0)
$insertStart = $find.Start

# The "_GoBack" bookmark currently marks that same spot; remove it now and
# re-create it at the new end of the document once all the new content has
# been inserted.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$ins = $d.Range($insertStart, $insertStart)
$xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml' w14:paraId='01F9D1E9' w14:textId='7B43928C' w:rsidR='00B62F66' w:rsidRDefault='009A2FF8' w:rsidP='00053B6C'><w:pPr><w:pStyle w:val='Listenabsatz'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr></w:pPr><w:r><w:t xml:space='preserve'>wie </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>margin</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t xml:space='preserve'> nur im Element… Sozusagen ein Rand im Element und nicht drum herum. Wird auf </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>width</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t xml:space='preserve'> drauf gerechnet</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='Listenabsatz'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr></w:pPr><w:r><w:t>Man sollte für Gemeinsamkeiten immer Klassen anlegen. Spart Code und ist eleganter bei Änderungen</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='Listenabsatz'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr></w:pPr><w:r><w:t xml:space='preserve'>Borders: </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>border</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t xml:space='preserve'>: 1px </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>black</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t xml:space='preserve'> solid;</w:t></w:r><w:r><w:t xml:space='preserve'> (Weite Farbe Style)</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='Listenabsatz'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr></w:pPr><w:r><w:t>Border-radius: abgerundete Ecken</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='Listenabsatz'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr></w:pPr></w:p>"
[void]$ins.InsertXML($xml)

# The first inserted paragraph above must not start a new bullet - it is a
# continuation of the existing "Padding: " paragraph. Deleting the
# paragraph mark between them merges the two paragraphs back into one
# (Word keeps the *second* paragraph's properties, which is why that
# fragment above was given the original paragraph's identity attributes).
$mergeRange = $d.Range($insertStart, $insertStart + 1)
$mergeRange.Delete()

# Re-create "_GoBack" at the new end of the document (inside the final,
# now-empty, bullet paragraph). A temporary trailing character is used
# while positioning the bookmark because collapsed ranges placed exactly
# at the document's last position get mis-resolved by Bookmarks.Add.
$endPos = $d.Content.End
$endRange = $d.Range($endPos - 1, $endPos - 1)
$endRange.InsertAfter("x")
$newEnd = $d.Content.End
$bmRange = $d.Range($newEnd - 2, $newEnd - 2)
$d.Bookmarks.Add("_GoBack", $bmRange)
$tmpRange = $d.Range($newEnd - 2, $newEnd - 1)
$tmpRange.Delete()
